$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "（此时所有页面只可能属于1，2类）",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "（此时被修改后页面只可能属于1，2类，之后可能会被访问或修改变成3，4类）",
    2
)
